$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.22405966666667
$ws.Range("N2").Value = 84.672179
$ws.Range("O2").Value = 0.3816548478108986
$ws.Range("P2").Value = 0.3816548478108986
$ws.Range("Q2").Value = 2.051625713209777
$ws.Range("R2").Value = 18.464631418888
$ws.Range("S2").Value = 0.3816548478108986
$ws.Range("T2").Value = 0.3816548478108986

# Row 3 updates
$ws.Range("N3").Value = 59.306181
$ws.Range("O3").Value = 0.2673191094302723
$ws.Range("P3").Value = 0.2673191094302723
$ws.Range("Q3").Value = 1.437001944781333
$ws.Range("S3").Value = 0.2673191094302723
$ws.Range("T3").Value = 0.2673191094302723

# Row 4 updates
$ws.Range("M4").Value = 25.95900466666667
$ws.Range("N4").Value = 77.877014
$ws.Range("O4").Value = 0.351026042758829
$ws.Range("P4").Value = 0.351026042758829
$ws.Range("S4").Value = 0.351026042758829
$ws.Range("T4").Value = 0.351026042758829
